$wb = $excel.ActiveWorkbook

# Helper: write a value that looks like a plain number but must be stored as
# TEXT (shared string), matching the original workbook's convention of
# keeping these numeric-looking results as strings. A leading apostrophe
# forces text entry; ClearFormats() then strips the transient "quote
# prefix" cell style that Excel applies, so no stray style survives.
function Set-TextValue($range, [string]$text) {
    $range.Formula = "'" + $text
    $range.ClearFormats()
}

# Sheet "Restricciones_del_follower": update Gamma_value row expressions/values
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
$wsFollower.Range("A2").Value = "-1 + 11.175706019321918y"
Set-TextValue $wsFollower.Range("B2") "35.87982986376233"
Set-TextValue $wsFollower.Range("E2") "0.7000000000000001"
Set-TextValue $wsFollower.Range("F2") "6.0"

# Row 3 (J_0_LP_v)
$wsFollower.Range("A3").Value = "-1 + 0.4631211675015009y"
Set-TextValue $wsFollower.Range("B3") "0.5282998527549532"
Set-TextValue $wsFollower.Range("E3") "1.4000000000000001"
Set-TextValue $wsFollower.Range("F3") "1.6"

# Sheet "Vector_bf" (sheet index 5). Note: worksheet names "Vector_bf" and
# "Vector_BF" differ only by case, and name-based lookup is case-insensitive,
# so we must use positional indices to disambiguate them.
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf.Range("A2") "-46.66302266159494"

# Sheet "Vector_BF" (sheet index 6)
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "1.4634002944900941"
Set-TextValue $wsBF.Range("A3") "-109.23755668724215"

# Sheet "Vector_Alpha" - this cell is a genuine numeric value (not text)
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 0.5368788324984991
